$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.62872295045842
$ws.Range("C2").Value = 10.69096673853132
$ws.Range("D2").Value = 6.27619033364414
$ws.Range("E2").Value = 12.04723991820016
$ws.Range("F2").Value = 45.50761799847239
$ws.Range("I2").Value = 30.69713944760976
$ws.Range("J2").Value = 10.22790188619332
$ws.Range("K2").Value = 14.81760398917551
$ws.Range("M2").Value = 17.71410498327276
$ws.Range("N2").Value = 23.10245998556939
$ws.Range("B3").Value = 13.45722975573752
$ws.Range("C3").Value = 10.57631360748561
$ws.Range("D3").Value = 6.269605433283227
$ws.Range("E3").Value = 12.05130322508275
$ws.Range("F3").Value = 45.48662517687411
$ws.Range("I3").Value = 30.73301094569045
$ws.Range("J3").Value = 10.24619272929466
$ws.Range("K3").Value = 14.71138967262536
$ws.Range("M3").Value = 17.68454041789056
$ws.Range("N3").Value = 23.15580052701276
$ws.Range("B4").Value = 13.35452903539459
$ws.Range("C4").Value = 10.50828343281554
$ws.Range("D4").Value = 6.266529836636598
$ws.Range("E4").Value = 12.05570197436559
$ws.Range("F4").Value = 45.48284905741044
$ws.Range("I4").Value = 30.76015662726947
$ws.Range("J4").Value = 10.25852432375581
$ws.Range("K4").Value = 14.64942813518405
$ws.Range("M4").Value = 17.66980029503447
$ws.Range("N4").Value = 23.19046458070316
$ws.Range("B5").Value = 13.31338244821356
$ws.Range("C5").Value = 10.4811867290635
$ws.Range("D5").Value = 6.265521279454934
$ws.Range("E5").Value = 12.05797392139188
$ws.Range("F5").Value = 45.4836034910418
$ws.Range("I5").Value = 30.77250464953056
$ws.Range("J5").Value = 10.26382665022762
$ws.Range("K5").Value = 14.62502032929413
$ws.Range("M5").Value = 17.664656379231
$ws.Range("N5").Value = 23.20507201218722
$ws.Range("B6").Value = 13.30659414247078
$ws.Range("C6").Value = 10.47672603467405
$ws.Range("D6").Value = 6.265368633587562
$ws.Range("E6").Value = 12.05838014925652
$ws.Range("F6").Value = 45.48386728991948
$ws.Range("I6").Value = 30.77463264334821
$ws.Range("J6").Value = 10.26472384128545
$ws.Range("K6").Value = 14.62101894965724
$ws.Range("M6").Value = 17.66385447089118
$ws.Range("N6").Value = 23.20752666645696
$ws.Range("B7").Value = 13.35397119707696
$ws.Range("C7").Value = 10.50791542220359
$ws.Range("D7").Value = 6.266515241942431
$ws.Range("E7").Value = 12.05573067268862
$ws.Range("F7").Value = 45.48284994575481
$ws.Range("I7").Value = 30.76031795274316
$ws.Range("J7").Value = 10.25859471047214
$ws.Range("K7").Value = 14.64909552353875
$ws.Range("M7").Value = 17.66972742325446
$ws.Range("N7").Value = 23.19065963113141
$ws.Range("B8").Value = 13.56908461594836
$ws.Range("C8").Value = 10.6509619242484
$ws.Range("D8").Value = 6.273719962790787
$ws.Range("E8").Value = 12.04824624827991
$ws.Range("F8").Value = 45.4984894546978
$ws.Range("I8").Value = 30.70844428241137
$ws.Range("J8").Value = 10.23398022626676
$ws.Range("K8").Value = 14.7803204673493
$ws.Range("M8").Value = 17.70320581625928
$ws.Range("N8").Value = 23.12045515821165
$ws.Range("B9").Value = 14.00919210843042
$ws.Range("C9").Value = 10.94884274186835
$ws.Range("D9").Value = 6.295456120891131
$ws.Range("E9").Value = 12.04864105365568
$ws.Range("F9").Value = 45.60135339612514
$ws.Range("I9").Value = 30.64741379164499
$ws.Range("J9").Value = 10.19443533245925
$ws.Range("K9").Value = 15.06236727845415
$ws.Range("M9").Value = 17.79570304893642
$ws.Range("N9").Value = 22.99793482200433
$ws.Range("B10").Value = 14.34050387788853
$ws.Range("C10").Value = 11.17631988751211
$ws.Range("D10").Value = 6.315963619154025
$ws.Range("E10").Value = 12.05806585061855
$ws.Range("F10").Value = 45.72071058970767
$ws.Range("I10").Value = 30.62745983616441
$ws.Range("J10").Value = 10.17068468597037
$ws.Range("K10").Value = 15.28308342396595
$ws.Range("M10").Value = 17.8796784461566
$ws.Range("N10").Value = 22.91711806616632
$ws.Range("B11").Value = 14.49227977336347
$ws.Range("C11").Value = 11.28124851244455
$ws.Range("D11").Value = 6.326254530029028
$ws.Range("E11").Value = 12.06432241227009
$ws.Range("F11").Value = 45.78443998119154
$ws.Range("I11").Value = 30.62379643886291
$ws.Range("J11").Value = 10.1610280284165
$ws.Range("K11").Value = 15.38606779861135
$ws.Range("M11").Value = 17.92127179712264
$ws.Range("N11").Value = 22.88234227467576
$ws.Range("B12").Value = 14.54985007968526
$ws.Range("C12").Value = 11.32115428197619
$ws.Range("D12").Value = 6.330287523231497
$ws.Range("E12").Value = 12.06697323346817
$ws.Range("F12").Value = 45.809920208102
$ws.Range("I12").Value = 30.62318792624927
$ws.Range("J12").Value = 10.15753604771355
$ws.Range("K12").Value = 15.4254052126034
$ws.Range("M12").Value = 17.93750154121341
$ws.Range("N12").Value = 22.86945898187848
$ws.Range("B13").Value = 14.53744789852855
$ws.Range("C13").Value = 11.31255281130716
$ws.Range("D13").Value = 6.329412932750493
$ws.Range("E13").Value = 12.06638983433504
$ws.Range("F13").Value = 45.80437282529568
$ws.Range("I13").Value = 30.62328434550841
$ws.Range("J13").Value = 10.15828078327838
$ws.Range("K13").Value = 15.41691860505739
$ws.Range("M13").Value = 17.93398500407323
$ws.Range("N13").Value = 22.87222093515787
$ws.Range("B14").Value = 14.49701451119306
$ws.Range("C14").Value = 11.28452835481797
$ws.Range("D14").Value = 6.326583613668935
$ws.Range("E14").Value = 12.06453486287509
$ws.Range("F14").Value = 45.7865093106699
$ws.Range("I14").Value = 30.62373076893133
$ws.Range("J14").Value = 10.16073743992181
$ws.Range("K14").Value = 15.38929748992905
$ws.Range("M14").Value = 17.92259747450541
$ws.Range("N14").Value = 22.88127663761526
$ws.Range("B15").Value = 14.47225874809993
$ws.Range("C15").Value = 11.26738382522015
$ws.Range("D15").Value = 6.324868221900106
$ws.Range("E15").Value = 12.06343526216086
$ws.Range("F15").Value = 45.77574255262002
$ws.Range("I15").Value = 30.62410563211809
$ws.Range("J15").Value = 10.16226366572359
$ws.Range("K15").Value = 15.37242202540731
$ws.Range("M15").Value = 17.9156844303806
$ws.Range("N15").Value = 22.88686069316153
$ws.Range("B16").Value = 14.33060145947363
$ws.Range("C16").Value = 11.16948855794543
$ws.Range("D16").Value = 6.31531024879347
$ws.Range("E16").Value = 12.0576964642312
$ws.Range("F16").Value = 45.71673478991231
$ws.Range("I16").Value = 30.62780821505654
$ws.Range("J16").Value = 10.17133884401626
$ws.Range("K16").Value = 15.27640238179129
$ws.Range("M16").Value = 17.87702776796241
$ws.Range("N16").Value = 22.91943072234547
$ws.Range("B17").Value = 14.24393046568112
$ws.Range("C17").Value = 11.10977776209567
$ws.Range("D17").Value = 6.309691537446062
$ws.Range("E17").Value = 12.0546791881008
$ws.Range("F17").Value = 45.68294540673254
$ws.Range("I17").Value = 30.6314665040785
$ws.Range("J17").Value = 10.1771999381047
$ws.Range("K17").Value = 15.21813590415558
$ws.Range("M17").Value = 17.85417613926878
$ws.Range("N17").Value = 22.93992033954447
$ws.Range("B18").Value = 14.19418274715001
$ws.Range("C18").Value = 11.07557218989582
$ws.Range("D18").Value = 6.30655049049154
$ws.Range("E18").Value = 12.05312918313305
$ws.Range("F18").Value = 45.66439918325084
$ws.Range("I18").Value = 30.6340802300175
$ws.Range("J18").Value = 10.18067911979567
$ws.Range("K18").Value = 15.18486728476959
$ws.Range("M18").Value = 17.8413523602118
$ws.Range("N18").Value = 22.95189257957915
$ws.Range("B19").Value = 14.1773584752545
$ws.Range("C19").Value = 11.06401567726393
$ws.Range("D19").Value = 6.30550262898682
$ws.Range("E19").Value = 12.0526362747421
$ws.Range("F19").Value = 45.65827261516255
$ws.Range("I19").Value = 30.6350526997736
$ws.Range("J19").Value = 10.18187567336285
$ws.Range("K19").Value = 15.17364608193938
$ws.Range("M19").Value = 17.83706563960062
$ws.Range("N19").Value = 22.9559783292834
$ws.Range("B20").Value = 14.25314647783613
$ws.Range("C20").Value = 11.11612001024473
$ws.Range("D20").Value = 6.31028028920127
$ws.Range("E20").Value = 12.05498119949227
$ws.Range("F20").Value = 45.68645044624298
$ws.Range("I20").Value = 30.63102433097563
$ws.Range("J20").Value = 10.17656483510008
$ws.Range("K20").Value = 15.22431336902183
$ws.Range("M20").Value = 17.85657568881551
$ws.Range("N20").Value = 22.93771981762217
$ws.Range("B21").Value = 14.50888861433462
$ws.Range("C21").Value = 11.29275545174447
$ws.Range("D21").Value = 6.327410979703176
$ws.Range("E21").Value = 12.06507208381921
$ws.Range("F21").Value = 45.79171977297896
$ws.Range("I21").Value = 30.62357850859191
$ws.Range("J21").Value = 10.16001139004115
$ws.Range("K21").Value = 15.39740152120511
$ws.Range("M21").Value = 17.92592933286056
$ws.Range("N21").Value = 22.87860901137208
$ws.Range("B22").Value = 14.67656549106364
$ws.Range("C22").Value = 11.40918017596763
$ws.Range("D22").Value = 6.339398730354253
$ws.Range("E22").Value = 12.07330757803516
$ws.Range("F22").Value = 45.86836716450355
$ws.Range("I22").Value = 30.62325117856182
$ws.Range("J22").Value = 10.15015313724381
$ws.Range("K22").Value = 15.51248839313588
$ws.Range("M22").Value = 17.97404495263572
$ws.Range("N22").Value = 22.84164090774323
$ws.Range("B23").Value = 14.58704262209082
$ws.Range("C23").Value = 11.34696419554624
$ws.Range("D23").Value = 6.332928971172897
$ws.Range("E23").Value = 12.06876259344231
$ws.Range("F23").Value = 45.82674432728811
$ws.Range("I23").Value = 30.62301057394379
$ws.Range("J23").Value = 10.15532687947685
$ws.Range("K23").Value = 15.45089519528945
$ws.Range("M23").Value = 17.94811253876589
$ws.Range("N23").Value = 22.8612193139642
$ws.Range("B24").Value = 14.24897966337065
$ws.Range("C24").Value = 11.11325229404644
$ws.Range("D24").Value = 6.310013836451228
$ws.Range("E24").Value = 12.05484408481493
$ws.Range("F24").Value = 45.68486307672531
$ws.Range("I24").Value = 30.63122264709631
$ws.Range("J24").Value = 10.17685162387355
$ws.Range("K24").Value = 15.2215198196991
$ws.Range("M24").Value = 17.85548987395003
$ws.Range("N24").Value = 22.93871407390384
$ws.Range("B25").Value = 13.88850200416767
$ws.Range("C25").Value = 10.86659881047857
$ws.Range("D25").Value = 6.288771260488072
$ws.Range("E25").Value = 12.04692558802497
$ws.Range("F25").Value = 45.56581523101211
$ws.Range("I25").Value = 30.65955843543721
$ws.Range("J25").Value = 10.2042008724488
$ws.Range("K25").Value = 14.98357797878301
$ws.Range("M25").Value = 17.76783943840669
$ws.Range("N25").Value = 23.0294617783635
